$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 342, shifting rows 342:485 down to 343:486
$ws.Rows("342:342").Insert()

# Populate the newly inserted row 342 with the new data
$ws.Range("A342").Value = 3
$ws.Range("B342").Value = "Femacal de La Calera"
$ws.Range("C342").Value = "Coquimbo"
$ws.Range("D342").Value = 44839
$ws.Range("E342").Value = 5
$ws.Range("F342").Value = "Fruta"
$ws.Range("G342").Value = 100108
$ws.Range("H342").Value = "Tropicales y subtropicales"
$ws.Range("I342").Value = 100108002
$ws.Range("J342").Value = "Mango"
$ws.Range("K342").Value = "Sin especificar"
$ws.Range("L342").Value = "Primera"
$ws.Range("M342").Value = 228
$ws.Range("N342").Value = 9000
$ws.Range("O342").Value = 9000
$ws.Range("P342").Value = 9000
$ws.Range("Q342").Value = "`$/bandeja 4 kilos"
$ws.Range("R342").Value = "Brasil"
$ws.Range("S342").Value = 2250
$ws.Range("T342").Value = 4
